$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "62.645.49"
$ws.Range("E2").Value = "  -1.99%  "

Set-TextValue "D3" "3.203.87"
$ws.Range("E3").Value = "  -3.03%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue "D5" "596.89"
$ws.Range("E5").Value = "  -0.32%  "

Set-TextValue "D6" "136.30"
$ws.Range("E6").Value = "  -4.67%  "

Set-TextValue "D7" "0.998"
$ws.Range("E7").Value = "  -0.32%  "

Set-TextValue "D8" "3.202.95"
$ws.Range("E8").Value = "  -2.94%  "

$ws.Range("E9").Value = "  -2.73%  "

$ws.Range("E10").Value = "  -3.17%  "

$ws.Range("E11").Value = "  -1.51%  "

$ws.Range("E12").Value = "  -3.18%  "

$ws.Range("E13").Value = "  -3.51%  "

Set-TextValue "D14" "33.70"
$ws.Range("E14").Value = "  -3.18%  "

Set-TextValue "D15" "3.731.14"

$ws.Range("E16").Value = "  -0.20%  "

Set-TextValue "D17" "3.222.14"
$ws.Range("E17").Value = "  -2.62%  "

Set-TextValue "D18" "62.768.97"
$ws.Range("E18").Value = "  -1.93%  "

$ws.Range("E19").Value = "  -2.43%  "

Set-TextValue "D20" "463.77"
$ws.Range("E20").Value = "  -3.59%  "

Set-TextValue "D21" "14.04"
$ws.Range("E21").Value = "  -1.66%  "

$ws.Range("E22").Value = "  -3.61%  "

$ws.Range("E23").Value = "  -3.82%  "

Set-TextValue "D24" "13.67"
$ws.Range("E24").Value = "  +1.39%  "

Set-TextValue "D25" "83.68"
$ws.Range("E25").Value = "  -0.78%  "

$ws.Range("E27").Value = "  -1.67%  "

$ws.Range("E28").Value = "  -0.11%  "

Set-TextValue "D29" "7.97"
$ws.Range("E29").Value = "  -3.34%  "

$ws.Range("E30").Value = "  -5.55%  "

Set-TextValue "D31" "2.09"
$ws.Range("E31").Value = "  -2.77%  "

Set-TextValue "D32" "27.63"
$ws.Range("E32").Value = "  -3.05%  "

$ws.Range("E33").Value = "  -3.69%  "

$ws.Range("E34").Value = "  -4.11%  "

$ws.Range("E35").Value = "  -4.58%  "

Set-TextValue "D36" "5.89"
$ws.Range("E36").Value = "  -1.71%  "

Set-TextValue "D37" "51.65"
$ws.Range("E37").Value = "  -3.17%  "

Set-TextValue "D38" "0.0₃0702"
$ws.Range("E38").Value = "  -4.81%  "

Set-TextValue "D39" "0.0394"
$ws.Range("E39").Value = "  -1.01%  "

Set-TextValue "D40" "421.36"
$ws.Range("E40").Value = "  -2.43%  "

Set-TextValue "D41" "3.020.64"
$ws.Range("E41").Value = "  +0.16%  "

$ws.Range("E42").Value = "  +5.68%  "

$ws.Range("E43").Value = "  -3.38%  "

$ws.Range("E44").Value = "  -3.92%  "

$ws.Range("E45").Value = "  -5.10%  "

Set-TextValue "D46" "2.18"
$ws.Range("E46").Value = "  -1.91%  "

Set-TextValue "D47" "36.38"
$ws.Range("E47").Value = "  +2.38%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D48" "26.17"
$ws.Range("E48").Value = "  -0.41%  "

$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D49" "0.999"
$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D50" "2.32"
$ws.Range("E50").Value = "  -0.17%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D51" "125.75"
$ws.Range("E51").Value = "  +1.31%  "
